# Course catalogue v2.6 - reintroduce licence metadata from repo (github api)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column L: license_name ------------------------------------------

$ws.Range("L1").Value = "license_name"
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108

$ws.Range("L2").Value  = "NULL"
$ws.Range("L3").Value  = "NULL"
$ws.Range("L4").Value  = "NULL"
$ws.Range("L5").Value  = "Other"
$ws.Range("L6").Value  = "NULL"
$ws.Range("L7").Value  = "NULL"
$ws.Range("L8").Value  = "NULL"
$ws.Range("L9").Value  = "MIT License"
$ws.Range("L10").Value = "NULL"
$ws.Range("L11").Value = "NULL"
$ws.Range("L12").Value = "NULL"
$ws.Range("L13").Value = "NULL"
$ws.Range("L14").Value = "NULL"

# --- Row 9: Fundemental theories in Machine Learning ----------------------
# course_objective / learning_objective_detail content reshuffled

$ws.Range("E9").Value = "Provide a thorough introduction into probability theory and statistical inference including maximum-likelihood and Bayesian approaches. Introduce supervised learning methods: linear and nonlinear regressions and classification algorithms. Introduce unsupervised learning methods: clustering, and dimensionality reduction. Brief introduction to Directed Graphical Models with a case study/example."

$ws.Range("F9").Value = "Be able to describe the diffence between frequentist and bayesian statistics; Understand the fundementals of probability theory, bayesian rule and inference, and the characteristics of major probability distributions.; Get a good understanding of major supervised learning algorithms specifically linear in parameter regression, bayesian linear regression, and classification methods.; Get a good understanding of main unsupervised learning algorithms specifically clustering and data dimensionality reduction algorithms.; Get familiar with directed graphical method as a technique of combining supervised and unsupervised learning into one modelling framework; be prepared to build on their current knowledge or take on more advanced courses such as application of machine learning techniques in natural language processing; be prepared to apply their knowledge through formulating machine learning problems and coding using standard libraries (e.g. in R and Python)"

# --- Row 11: Introduction to Reproducibility ------------------------------

$ws.Range("C11").Value = "6 Hours"
$ws.Range("D11").Value = 'Introduction to reproducibility. What is a reproducible analytical pipeline (RAP), writing "good code" and creating a reproducible report.'
$ws.Range("E11").Value = "Participants should gain an awareness of the importance of reproducibility in their work. Learners will also gain experience of linting code in Python and using parameterised reports in R markdown."
$ws.Range("F11").Value = "represent pipelines and identify opportunities to automate; consider adherance to a programming style guide; use linting software to standardise Python scripts; use parameterised R markdown reports to improve the efficiency of report production."
$ws.Range("H11").Value = "Basic familiarity with Python and R programming is assumed."
